# example and document for range in sheet #1811
# Add a new sheet "PatientsShifted" that mirrors the "Patients" sheet's
# data/styles but shifted right by one column and down by five rows, and
# update the selections / active-sheet state on the existing sheets.

$wb = $excel.ActiveWorkbook

$wsPatients  = $wb.Worksheets.Item("Patients")
$wsProviders = $wb.Worksheets.Item("Providers")

# --- Update the (soon to be non-active) sheets' selections first -----------
# Patients: selection moves to B2, and it stops being the tab-selected sheet
# once PatientsShifted is activated below.
$wsPatients.Range("B2").Select() | Out-Null

# Providers: selection moves to F57.
$wsProviders.Range("F57").Select() | Out-Null

# --- Add the new sheet at the end ------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsShifted = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsShifted.Name = "PatientsShifted"

# --- Copy the Patients data into PatientsShifted, shifted +1 col / +5 rows -
# Source used range on Patients is A1:O5; destination top-left is B6.
$srcRange = $wsPatients.Range("A1:O5")
$dstRange = $wsShifted.Range("B6")

# Paste formats first (so the later value paste doesn't overwrite them with
# an auto-generated date number format), then paste the values.
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)   # xlPasteFormats

$srcRange.Copy()
$dstRange.PasteSpecial(-4163)   # xlPasteValues

$excel.CutCopyMode = $false

# A few source cells are blank-but-unstyled; the block paste above leaves
# behind empty styled placeholders at the equivalent shifted positions.
# Clear them so the sheet matches the source's sparse cell layout exactly.
$wsShifted.Range("N8").Clear() | Out-Null
$wsShifted.Range("G9").Clear() | Out-Null
$wsShifted.Range("O9").Clear() | Out-Null

# --- Make PatientsShifted the active sheet/selection ------------------------
$wsShifted.Activate()
$wsShifted.Range("G7").Select() | Out-Null
